# Dataframe ST.xlsx — add a new "21-nov" snapshot column (CM) to Sheet1 and
# refresh the VLOOKUP source table on Sheet3 (rows 20:36) with updated figures.
# Sheet1's CB/CC columns carry live VLOOKUP formulas against Sheet3 and will
# recalc automatically; CL (an earlier frozen "17-nov" snapshot) is left as-is,
# and the new CM column is populated with the freshly computed figures.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- 1. Update the lookup table on Sheet3 (rows 20..36, column B) ----------
$ws3.Range("B20").Value = 7.8124694268455457
$ws3.Range("B21").Value = 1.0789137669908284
$ws3.Range("B22").Value = 7.1933541229111588
$ws3.Range("B23").Value = 8.241300995293086
$ws3.Range("B24").Value = 7.4503895899621275
$ws3.Range("B25").Value = 0
$ws3.Range("B26").Value = 13.093917029548642
$ws3.Range("B27").Value = 12.247737578050847
$ws3.Range("B28").Value = 4.8564437573668622
$ws3.Range("B29").Value = 0.84188632979554467
$ws3.Range("B30").Value = 12.258000000037805
$ws3.Range("B31").Value = 12.477850000103997
$ws3.Range("B32").Value = 12.709138617224118
$ws3.Range("B33").Value = 10.864299891629882
$ws3.Range("B34").Value = 12.852497738105013
$ws3.Range("B35").Value = 5.7974960181187329
$ws3.Range("B36").Value = 34.878552096152589

# --- 2. Add the new "21-nov" column header on Sheet1 (CM1) -----------------
$ws1.Range("CM1").Value = "21-nov"

# --- 3. Populate the new CM column with the refreshed figures --------------
$ws1.Range("CM2").Value = 4.8564437573668622
$ws1.Range("CM3").Value = 0.84188632979554467
$ws1.Range("CM4").Value = 13.093917029548642
$ws1.Range("CM5").Value = 12.247737578050847
$ws1.Range("CM6").Value = 8.241300995293086
$ws1.Range("CM7").Value = 7.4503895899621275
$ws1.Range("CM8").Value = 7.1933541229111588
$ws1.Range("CM9").Value = 5.7974960181187329
$ws1.Range("CM10").Value = 10.864299891629882
$ws1.Range("CM11").Value = 12.852497738105013
$ws1.Range("CM12").Value = 7.8124694268455457
$ws1.Range("CM13").Value = 1.0789137669908284
$ws1.Range("CM14").Value = 34.878552096152589
$ws1.Range("CM15").Value = 12.709138617224118
$ws1.Range("CM16").Value = 0
$ws1.Range("CM17").Value = 12.258000000037805
$ws1.Range("CM18").Value = 12.477850000103997

# Match the numeric format of the existing CL column so the style reuses the
# same cell-format record (numFmtId 1, the integer-style already used by CL).
$ws1.Range("CM2:CM18").NumberFormat = "0"

# --- 4. Leave the selection where the author left it after the edit --------
$ws1.Range("CM2").Select()
